# "addind preparer to sheet"
# The "purpose" column (E2:E16) is updated from "S.GISH" to a new value
# "fullRNASEQ" for every data row, and the sheet's visible selection /
# scroll position is moved down to the newly-relevant area (D17:F24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "purpose" column for every data row (2-16) to the new value.
$ws.Range("E2:E16").Value = "fullRNASEQ"

# Scroll the view down a bit (best effort - mirrors the author's saved
# view position before re-selecting a new range).
try {
    $excel.ActiveWindow.ScrollRow = 3
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

# Move the active selection to D17:F24, with D17 as the active cell.
$null = $ws.Range("D17:F24").Select()
